$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q1" sheet right before "总计" ---
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Header row (same layout/style as the other quarterly sheets)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the header style from an existing quarterly sheet so formatting matches
$srcHeader = $wb.Worksheets.Item("2021-Q4").Range("B1:H1")
$srcHeader.Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$srcA2 = $wb.Worksheets.Item("2021-Q4").Range("A2")
$srcA2.Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "006195"
$newSheet.Range("C2").Value = "国金量化多因子股票"
$newSheet.Range("D2").Value = "0.09"
$newSheet.Range("E2").Value = "80.71"
$newSheet.Range("F2").Value = "0.88"
$newSheet.Range("G2").Value = "0.0008"
$newSheet.Range("H2").Value = 9

# --- 2. Update the "总计" (summary) sheet: insert a new first data row for 2022-Q1 ---
$totalSheet = $wb.Worksheets.Item("总计")

# Shift existing data rows (2021-Q4, 2021-Q3, 2021-Q2) down by one row.
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 8
$totalSheet.Range("D5").Value = 0.23

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 8
$totalSheet.Range("D4").Value = 0.91

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 7
$totalSheet.Range("D3").Value = 0.91

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0
